$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Vasp -> Cxcr2, ECs -> ECs)
$ws.Range("G2").Value = 10.37416333333333
$ws.Range("H2").Value = 31.12249
$ws.Range("I2").Value = 0.2339167445675379
$ws.Range("J2").Value = 0.2339167445675379
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.03970866666666666
$ws.Range("N2").Value = 0.119126
$ws.Range("Q2").Value = 0.4119441937488889
$ws.Range("R2").Value = 3.70749774374
$ws.Range("S2").Value = 0.2339167445675379
$ws.Range("T2").Value = 0.2339167445675379

# Row 3 (Vasp -> Cxcr2, FAPs -> ECs)
$ws.Range("I3").Value = 0.2155315552192661
$ws.Range("J3").Value = 0.2155315552192661
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.03970866666666666
$ws.Range("N3").Value = 0.119126
$ws.Range("Q3").Value = 0.3795665543584443
$ws.Range("R3").Value = 3.416098989226
$ws.Range("S3").Value = 0.2155315552192661
$ws.Range("T3").Value = 0.2155315552192661

# Row 4 (Vasp -> Cxcr2, MuSCs -> ECs)
$ws.Range("G4").Value = 24.41686366666667
$ws.Range("H4").Value = 73.250591
$ws.Range("I4").Value = 0.5505517002131959
$ws.Range("J4").Value = 0.550551700213196
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.03970866666666666
$ws.Range("N4").Value = 0.119126
$ws.Range("Q4").Value = 0.9695611003851111
$ws.Range("R4").Value = 8.726049903466
$ws.Range("S4").Value = 0.5505517002131959
$ws.Range("T4").Value = 0.550551700213196
